# update data on Feb-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Case 6 (row 7): discharged on Feb-20 ---
$ws.Range("K7").Value = "Discharged"
$ws.Range("N7").Value = "Feb-20"

# --- Case 16 (row 17): discharged on Feb-20 ---
$ws.Range("K17").Value = "Discharged"
$ws.Range("N17").Value = "Feb-20"

# --- Occupation data added for several existing cases ---
$ws.Range("P20").Value = "Sales Person"
$ws.Range("P21").Value = "Sales Person"
$ws.Range("P22").Value = "Maid"
$ws.Range("P26").Value = "Jewellry Clerk"
$ws.Range("P33").Value = "Teacher"
$ws.Range("P35").Value = "Sales Person"
$ws.Range("P36").Value = "Taxi Driver"
$ws.Range("P38").Value = "Private-Hire Driver"
$ws.Range("P45").Value = "Security Officer"

# --- Case 66 (row 67): discharged on Feb-20 ---
$ws.Range("K67").Value = "Discharged"
$ws.Range("N67").Value = "Feb-20"

$ws.Range("P79").Value = "Hospital administration"

# --- Case 83 (row 84): updated location/visited/symptom date ---
$ws.Range("B84").Value = 1.387508
$ws.Range("C84").Value = 103.90475499999999
$ws.Range("H84").Value = "Rivervale Drive"
$ws.Range("I84").Value = "Malaysia, Philemon Singapore Pte Ltd (16 Kallang Place), GP clinic"
$ws.Range("M84").Value = "Jan-28"

# --- Case 84 (row 85): updated location/visited/symptom date ---
$ws.Range("B85").Value = 1.323011
$ws.Range("C85").Value = 103.881266
$ws.Range("H85").Value = "Aljunied Road"
$ws.Range("I85").Value = "Lonza Biologics (35 Tuas South Avenue 6), Bugis Junction, GP clinic"
$ws.Range("M85").Value = "Feb-04"

# --- New case 85 (row 86) ---
$ws.Range("A86").Value = 85
$ws.Range("B86").Value = 1.322109
$ws.Range("C86").Value = 103.847272
$ws.Range("D86").Value = "Feb-20"
$ws.Range("E86").Value = 36
$ws.Range("F86").Value = "Male"
$ws.Range("G86").Value = "Singapore"
$ws.Range("K86").Value = "NCID"
$ws.Range("L86").Value = "Chinese"
$ws.Range("P86").Value = "Singapore Work Pass holder"

# --- Update view: scrolled down towards the newly-added row, C86 selected ---
$ws.Range("A77").Select()
$ws.Range("C86").Select()
